$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45205 -> 45206) for every data row (rows 2 through 537).
$ws.Range("C2:C537").Value = 45206
